$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$text = @'
questions = [
    {
        "title": "You are developing a game and need to optimize the way you write camera movement logic and physics-related calculations. Which of the following Update functions should you use?",
        "ques_type": 2,
        "options": [
            "LateUpdate for camera FixedUpdate for physics",
            "Update for camera FixedUpdate for physics",
            "FixedUpdate for camera LateUpdate for physics",
            "LateUpdate for camera Update for physics"
        ],
        "score": "LateUpdate for camera FixedUpdate for physics"
    },
    {
        "title": "In the below code block, what does yield return new WaitForSeconds(4.0f) do? void Start() {\n StartCoroutine(Up()) \n} IEnumerator Up() \n{ while (true) { \nyield return new WaitForSeconds(4.0f) \ntransform.Translate(0.0f, 10.0f, 0.0f)\n }\n }",
        "ques_type": 2,
        "options": [
            "It stops the coroutine after four seconds.",
            "It stops the coroutine immediately.",
            "It suspends the coroutine for four seconds.",
            "It stops the coroutine after four frames."
        ],
        "score": "It suspends the coroutine for four seconds."
    },
    {
        "title": "You are developing a multiplayer game. You need levels to be randomly generated and for all players in the game to be in the same environment. Which of the following should you do to achieve this?",
        "ques_type": 2,
        "options": [
            "Use the Random.Range((float min, float max)) function.",
            "Use both Random.Range and InitState functions.",
            "Use a flat file to store all of the level values.",
            "Use the InitState(int seed) function."
        ],
        "score": "Use the InitState(int seed) function."
    },
    {
        "title": "True or false: Though capable of making 2D and 2.5D games, the Unity engine is, at its core, a 3D engine.",
        "ques_type": 11,
        "options": [
            "true",
            "false"
        ],
        "score": "True"
    }
]
'@
$text = $text.TrimEnd("`r", "`n")

$ws.Range("A2").ClearContents()
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = $text
$ws.Rows(1).AutoFit()
